# Update the cryptos worksheet with refreshed price / volume(1h) data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows where only Price (D) and Volume 1h (E) change ---
# Price strings that look like plain decimal numbers (e.g. "1.001", "286.00")
# get a leading apostrophe so Excel stores them as text, exactly like the
# multi-dot prices (e.g. "22.393.49") are already stored as text.
$deChanges = @{
    2  = @("22.393.49",  "  -0.31%  ")
    3  = @("1.560.07",   "  -0.87%  ")
    4  = @("'1.001",     "  -0.08%  ")
    5  = @("'1.001",     "  +0.00%  ")
    6  = @("'286.00",    "  -1.86%  ")
    7  = @("'0.3642",    "  -2.85%  ")
    8  = @("'48.62",     "  -2.74%  ")
    9  = @("'0.3336",    "  -2.10%  ")
    11 = @("'0.07375",   "  -2.77%  ")
    12 = @("'1.001",     "  -0.05%  ")
    13 = @("'20.70",     "  -3.38%  ")
    14 = @("'5.891",     "  -1.83%  ")
    15 = @("'6.844",     "  -1.79%  ")
    16 = @("1.561.83",   "  -0.42%  ")
    17 = @("'0.00001097","  -2.53%  ")
    18 = @("'88.65",     "  -2.84%  ")
    19 = @("'0.06739",   "  -0.03%  ")
    20 = @("'1.001",     "  -0.04%  ")
    21 = @("'6.281",     "  -0.03%  ")
    22 = @("'15.97",     "  -2.89%  ")
    23 = @("'11.90",     "  -2.34%  ")
    24 = @("22.391.78",  "  -0.34%  ")
    25 = @("'2.391",     "  +2.69%  ")
    26 = @("'2.565",     "  -1.56%  ")
    27 = @("'149.41",    "  +0.40%  ")
    28 = @("'19.26",     "  -4.52%  ")
    29 = @("'5.009",     "  -0.01%  ")
    30 = @("'122.58",    "  -2.83%  ")
    31 = @("1.736.57",   "  -0.54%  ")
    32 = @("'1.059",     "  +1.52%  ")
    33 = @("'6.090",     "  -1.09%  ")
    35 = @("'9.548",     "  -3.66%  ")
    36 = @("'0.08218",   "  -2.67%  ")
    41 = @("'5.310",     "  -3.64%  ")
    42 = @("'11.07",     "  -2.96%  ")
    43 = @("'1.001",     "  +0.00%  ")
    44 = @("'0.6033",    "  -4.27%  ")
    45 = @("'13.59",     "  -4.12%  ")
    46 = @("'3.758",     "  -1.60%  ")
    47 = @("'0.5712",    "  -3.00%  ")
    48 = @("'124.52",    "  -4.51%  ")
    49 = @("'1.999",     "  -4.79%  ")
    50 = @("'1.203",     "  -2.02%  ")
    51 = @("'0.07213",   "  -1.81%  ")
}

foreach ($row in $deChanges.Keys) {
    $vals = $deChanges[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# --- Rows where only Volume (1h) changes ---
$eOnlyChanges = @{
    10 = "  -2.34%  "
    34 = "  +0.62%  "
}

foreach ($row in $eOnlyChanges.Keys) {
    $ws.Cells.Item($row, 5).Value = $eOnlyChanges[$row]
}

# --- Rows where Coin name/link/price/volume all change (coins re-ranked) ---
$fullChanges = @{
    37 = @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "'1.305",   "  -5.97%  ")
    38 = @("VeChain",          "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet",        "'0.02370", "  -4.13%  ")
    39 = @("Algorand",         "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo",       "'0.2208",  "  -3.99%  ")
    40 = @("Hedera",           "https://coinranking.com/coin/jad286TjB+hedera-hbar",             "'0.06351", "  -3.20%  ")
}

foreach ($row in $fullChanges.Keys) {
    $vals = $fullChanges[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
